# Weekly fruit/vegetable price update: a new "Primera"/"Segunda" price pair
# (week of 2021-10-19) is inserted right after the existing 2021-09-30 pair
# (rows 101-102), pushing every subsequent weekly pair down by one slot
# (2 rows). The two rows that fall off the bottom of the table become new
# rows 163-164 at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 103..162 down to 105..164, working from the bottom up so a
# source row is never clobbered before it has been copied.
for ($r = 162; $r -ge 103; $r--) {
    $dest = $r + 2
    $srcRange = $ws.Range("A" + $r + ":R" + $r)
    $dstRange = $ws.Range("A" + $dest + ":R" + $dest)
    $srcRange.Copy($dstRange)
}

# The freed-up rows 103/104 get a brand-new weekly pair, cloned from the
# previous week's pair (rows 101/102) and re-dated.
$ws.Range("A101:R101").Copy($ws.Range("A103:R103"))
$ws.Range("A102:R102").Copy($ws.Range("A104:R104"))

$ws.Range("D103").Value = 44488
$ws.Range("D104").Value = 44488
